$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 61: start 19:51 (0.82708333333333339), end 20:15 (0.84375)
$ws.Cells.Item(61, 4).Value = 0.82708333333333339
$ws.Cells.Item(61, 5).Value = 0.84375
$ws.Cells.Item(61, 7).Value = "Discovered bug in LogicGate.vhd that does not allow it to compile on ModelSim. Fixed it together with teammates."

# Row 62: start 20:15 (0.84375), end 20:39 (0.86041666666666661)
$ws.Cells.Item(62, 4).Value = 0.84375
$ws.Cells.Item(62, 5).Value = 0.86041666666666661
$ws.Cells.Item(62, 7).Value = "Verified project download works on a fresh installation of the project. Teammates had discovered this issue earlier and is a source of a massive headache."

# Row 63: start 20:39 (0.86041666666666661), end 21:10 (0.88194444444444453)
$ws.Cells.Item(63, 4).Value = 0.86041666666666661
$ws.Cells.Item(63, 5).Value = 0.88194444444444453
$ws.Cells.Item(63, 7).Value = "Re-compiled all VHD files and checked that all results are the same. Updated transcript , summary, .vho and .sdo files. DONE"

# Update the active selection to reflect the new working cell
$ws.Range("G64").Select()
